$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header row (row 1) labels for columns C through K so they
# match the data actually stored in each column.
$ws.Range("C1").Value = "nome"
$ws.Range("D1").Value = "logradouro"
$ws.Range("E1").Value = "numero"
$ws.Range("F1").Value = "complemento"
$ws.Range("G1").Value = "bairro"
$ws.Range("H1").Value = "cep"
$ws.Range("I1").Value = "cidade"
$ws.Range("J1").Value = "estado"
$ws.Range("K1").Value = "telefone"
